$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (interested count) column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 251
$wsExpo.Range("F4").Value = 881
$wsExpo.Range("F6").Value = 40

# Sheet "全部类型" (All types) - same rows shifted by one due to an extra
# performance entry present in this combined sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 251
$wsAll.Range("F5").Value = 881
$wsAll.Range("F7").Value = 40
